$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the IF formula for F4, and a shared formula F5:F10
$ws.Range("F4").Formula = '=IF(D4<E4,"BORÇLUYUZ","ALACAKLIYIZ")'
$ws.Range("F5:F10").Formula = '=IF(D5<E5,"BORÇLUYUZ","ALACAKLIYIZ")'

# F10 picks up the same cell style as F4:F9 (not the totals-row style)
$ws.Range("F9").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Student info in J9:J11
$ws.Range("J9").Value = 20215070019
$ws.Range("J10").Value = "KÜBRA ÇABUK"
$ws.Range("J11").Value = "YBS"

# Column F width adjustment (target stored width 16.28515625 characters;
# the engine quantizes to an MDW-7 pixel grid, so 15.57 is the closest
# settable value that rounds to the same stored width as native Excel)
$ws.Columns.Item(6).ColumnWidth = 15.57

# Update selection to H10
$ws.Range("H10").Select()
